$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric, so Excel
# stores them as text (matching the original inline-string cell type)
# instead of auto-converting to a number.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply the updated cell values.
$ws.Range('D2').Value = '43.760.28'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '2.280.93'
$ws.Range('E3').Value = '  +2.98%  '
$ws.Range('D5').Value = '265.66'
$ws.Range('E5').Value = '  +1.60%  '
$ws.Range('D6').Value = '93.30'
$ws.Range('E6').Value = '  +7.43%  '
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '0.618'
$ws.Range('E9').Value = '  +2.79%  '
$ws.Range('D10').Value = '45.97'
$ws.Range('E10').Value = '  +1.80%  '
$ws.Range('D11').Value = '0.0931'
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').Value = '7.98'
$ws.Range('E12').Value = '  +7.03%  '
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('D14').Value = '2.625.67'
$ws.Range('E14').Value = '  +3.11%  '
$ws.Range('D15').Value = '15.26'
$ws.Range('E15').Value = '  +5.49%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '0.832'
$ws.Range('E16').Value = '  +6.09%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.288.87'
$ws.Range('E17').Value = '  +3.60%  '
$ws.Range('D18').Value = '43.782.73'
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('D19').Value = '0.0000105'
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('D20').Value = '6.18'
$ws.Range('E20').Value = '  +3.98%  '
$ws.Range('D21').Value = '70.65'
$ws.Range('E21').Value = '  +1.16%  '
$ws.Range('B22').Value = 'ImmutableX'
$ws.Range('C22').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D22').Value = '2.27'
$ws.Range('E22').Value = '  -3.11%  '
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').Value = '10.00'
$ws.Range('E23').Value = '  +12.37%  '
$ws.Range('D24').Value = '233.69'
$ws.Range('E24').Value = '  +1.02%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').Value = '11.24'
$ws.Range('E26').Value = '  +5.53%  '
$ws.Range('E27').Value = '  +8.60%  '
$ws.Range('D28').Value = '3.37'
$ws.Range('E28').Value = '  -5.34%  '
$ws.Range('D29').Value = '39.51'
$ws.Range('E29').Value = '  -0.89%  '
$ws.Range('E30').Value = '  +1.00%  '
$ws.Range('D31').Value = '22.18'
$ws.Range('E31').Value = '  +8.34%  '
$ws.Range('D32').Value = '171.76'
$ws.Range('E32').Value = '  -1.76%  '
$ws.Range('D33').Value = '0.0912'
$ws.Range('E33').Value = '  +4.01%  '
$ws.Range('D34').Value = '5.54'
$ws.Range('E34').Value = '  +2.35%  '
$ws.Range('E35').Value = '  +0.93%  '
$ws.Range('E36').Value = '  -0.91%  '
$ws.Range('D37').Value = '4.41'
$ws.Range('E37').Value = '  -1.85%  '
$ws.Range('D38').Value = '0.0345'
$ws.Range('E38').Value = '  -3.43%  '
$ws.Range('E39').Value = '  +13.52%  '
$ws.Range('D40').Value = '0.242'
$ws.Range('E40').Value = '  +20.71%  '
$ws.Range('D41').Value = '2.30'
$ws.Range('E41').Value = '  +10.02%  '
$ws.Range('D42').Value = '12.28'
$ws.Range('E42').Value = '  -2.39%  '
$ws.Range('D43').Value = '1.30'
$ws.Range('E43').Value = '  +15.57%  '
$ws.Range('E44').Value = '  -2.14%  '
$ws.Range('D45').Value = '60.71'
$ws.Range('E45').Value = '  -4.96%  '
$ws.Range('D46').Value = '8.76'
$ws.Range('E46').Value = '  +5.36%  '
$ws.Range('E47').Value = '  +3.91%  '
$ws.Range('D48').Value = '99.21'
$ws.Range('E48').Value = '  -1.44%  '
$ws.Range('D49').Value = '1.18'
$ws.Range('E49').Value = '  -0.67%  '
$ws.Range('D50').Value = '2.504.14'
$ws.Range('E50').Value = '  +3.08%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = '0.420'
$ws.Range('E51').Value = '  -7.01%  '
